$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 196.18182
$ws.Range("I11").Value = 196.18182
$ws.Range("K11").Value = 196.18182
$ws.Range("M11").Value = -56.18181999999999
$ws.Range("H33").Value = 110.53846
$ws.Range("J33").Value = 111.85714
$ws.Range("L33").Value = 111.85714
$ws.Range("N33").Value = -569.85714
$ws.Range("H98").Value = 1607.4546
$ws.Range("I98").Value = 1636.381
$ws.Range("K98").Value = 1636.381
$ws.Range("M98").Value = -138.3810000000001
$ws.Range("H116").Value = 12780.23
$ws.Range("I116").Value = 16591.857
$ws.Range("K116").Value = 16591.857
$ws.Range("M116").Value = -13149.857
$ws.Range("H122").Value = 1607.4546
$ws.Range("I122").Value = 1636.381
$ws.Range("K122").Value = 4909.143
$ws.Range("M122").Value = -2459.143
$ws.Range("H129").Value = 909.65955
$ws.Range("I129").Value = 1500
$ws.Range("J129").Value = 883.42224
$ws.Range("K129").Value = 4500
$ws.Range("L129").Value = 2650.26672
$ws.Range("M129").Value = 500
$ws.Range("N129").Value = -12650.26672
$ws.Range("H132").Value = 1118.579
$ws.Range("I132").Value = 1045.8788
$ws.Range("K132").Value = 3137.6364
$ws.Range("M132").Value = -607.6363999999999
$ws.Range("H138").Value = 2758.1133
$ws.Range("J138").Value = 2827.6296
$ws.Range("L138").Value = 8482.8888
$ws.Range("N138").Value = -18762.8888
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3614.717
$ws.Range("I32").Value = 2532.5435
$ws.Range("J32").Value = 10726.143
$ws.Range("K32").Value = 2532.5435
$ws.Range("L32").Value = 10726.143
$ws.Range("M32").Value = -2245.5435
$ws.Range("N32").Value = -11300.143
$ws.Range("H37").Value = 14028.5
$ws.Range("J37").Value = 14028.5
$ws.Range("L37").Value = 14028.5
$ws.Range("N37").Value = -14574.5
$ws.Range("H61").Value = 3618.3225
$ws.Range("I61").Value = 2683.3572
$ws.Range("J61").Value = 12344.667
$ws.Range("K61").Value = 2683.3572
$ws.Range("L61").Value = 12344.667
$ws.Range("M61").Value = -2471.3572
$ws.Range("N61").Value = -12768.667
$ws.Range("H88").Value = 3819.5
$ws.Range("J88").Value = 4599.7144
$ws.Range("L88").Value = 4599.7144
$ws.Range("N88").Value = -5411.7144
$ws.Range("H91").Value = 3819.5
$ws.Range("J91").Value = 4599.7144
$ws.Range("L91").Value = 4599.7144
$ws.Range("N91").Value = -7407.7144
$ws.Range("H122").Value = 7997
$ws.Range("I122").Value = 7997
$ws.Range("K122").Value = 23991
$ws.Range("M122").Value = -21541
$ws.Range("H136").Value = 3618.3225
$ws.Range("I136").Value = 2683.3572
$ws.Range("J136").Value = 12344.667
$ws.Range("K136").Value = 8050.071599999999
$ws.Range("L136").Value = 37034.001
$ws.Range("M136").Value = -5500.071599999999
$ws.Range("N136").Value = -42134.001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 11867
$ws.Range("J75").Value = 11867
$ws.Range("L75").Value = 11867
$ws.Range("N75").Value = -13739
$ws.Range("H78").Value = 11867
$ws.Range("J78").Value = 11867
$ws.Range("L78").Value = 35601
$ws.Range("N78").Value = -44961
$ws.Range("H86").Value = 113454.28
$ws.Range("I86").Value = 2223.1667
$ws.Range("K86").Value = 2223.1667
$ws.Range("M86").Value = -1100.1667
$ws.Range("H89").Value = 113454.28
$ws.Range("I89").Value = 2223.1667
$ws.Range("K89").Value = 11115.8335
$ws.Range("M89").Value = -5499.833500000001
$ws.Range("H94").Value = 777.2778
$ws.Range("I94").Value = 766.5333
$ws.Range("K94").Value = 766.5333
$ws.Range("M94").Value = -315.5333000000001
$ws.Range("H134").Value = 11961
$ws.Range("I134").Value = 12117.3
$ws.Range("K134").Value = 36351.89999999999
$ws.Range("M134").Value = -33816.89999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2941.4412
$ws.Range("I31").Value = 1997.0476
$ws.Range("K31").Value = 1997.0476
$ws.Range("M31").Value = -1702.0476
$ws.Range("H34").Value = 2941.4412
$ws.Range("I34").Value = 1997.0476
$ws.Range("K34").Value = 1997.0476
$ws.Range("M34").Value = -1795.0476
$ws.Range("H99").Value = 1431130.4
$ws.Range("I99").Value = 5000999.5
$ws.Range("J99").Value = 3182.8
$ws.Range("K99").Value = 5000999.5
$ws.Range("L99").Value = 3182.8
$ws.Range("M99").Value = -4999501.5
$ws.Range("N99").Value = -6178.8
$ws.Range("H107").Value = 462.05884
$ws.Range("I107").Value = 411.92307
$ws.Range("J107").Value = 625
$ws.Range("K107").Value = 411.92307
$ws.Range("L107").Value = 625
$ws.Range("M107").Value = 1508.07693
$ws.Range("N107").Value = -4465
$ws.Range("H126").Value = 1431130.4
$ws.Range("I126").Value = 5000999.5
$ws.Range("J126").Value = 3182.8
$ws.Range("K126").Value = 15002998.5
$ws.Range("L126").Value = 9548.400000000001
$ws.Range("M126").Value = -15000528.5
$ws.Range("N126").Value = -14488.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 921.875
$ws.Range("J5").Value = 925
$ws.Range("L5").Value = 2775
$ws.Range("N5").Value = -2999
$ws.Range("H26").Value = 661
$ws.Range("I26").Value = 1150
$ws.Range("J26").Value = 335
$ws.Range("K26").Value = 3450
$ws.Range("L26").Value = 1005
$ws.Range("M26").Value = -3162
$ws.Range("N26").Value = -1581
$ws.Range("H132").Value = 500
$ws.Range("I132").Value = 500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970
$ws.Range("H135").Value = 921.875
$ws.Range("J135").Value = 925
$ws.Range("L135").Value = 8325
$ws.Range("N135").Value = -13395
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1811.1538
$ws.Range("I22").Value = 2650
$ws.Range("J22").Value = 1438.3334
$ws.Range("K22").Value = 2650
$ws.Range("L22").Value = 1438.3334
$ws.Range("M22").Value = -2355
$ws.Range("N22").Value = -2028.3334
$ws.Range("H27").Value = 1811.1538
$ws.Range("I27").Value = 2650
$ws.Range("J27").Value = 1438.3334
$ws.Range("K27").Value = 2650
$ws.Range("L27").Value = 1438.3334
$ws.Range("M27").Value = -2543
$ws.Range("N27").Value = -1652.3334
$ws.Range("H55").Value = 273.34286
$ws.Range("I55").Value = 240.2
$ws.Range("K55").Value = 240.2
$ws.Range("M55").Value = -67.19999999999999
$ws.Range("H68").Value = 2498.2307
$ws.Range("I68").Value = 2225.182
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2225.182
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1476.182
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 2498.2307
$ws.Range("I71").Value = 2225.182
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 11125.91
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -7381.91
$ws.Range("N71").Value = -27488
$ws.Range("H93").Value = 489.69232
$ws.Range("I93").Value = 416.6
$ws.Range("K93").Value = 416.6
$ws.Range("M93").Value = 831.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 58840972
$ws.Range("J2").Value = 18531.812
$ws.Range("L2").Value = 18531.812
$ws.Range("N2").Value = -18755.812
$ws.Range("H122").Value = 256935.2
$ws.Range("I122").Value = 256935.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 770805.6000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -768355.6000000001
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2648.1482
$ws.Range("I132").Value = 2240.125
$ws.Range("K132").Value = 6720.375
$ws.Range("M132").Value = -4190.375
$ws.Range("H136").Value = 19159596
$ws.Range("I136").Value = 30866610
$ws.Range("J136").Value = 2666
$ws.Range("K136").Value = 92599830
$ws.Range("L136").Value = 7998
$ws.Range("M136").Value = -92597280
$ws.Range("N136").Value = -13098
